$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update expiry dates for FE008 (row 9) and FE009 (row 10)
$ws.Range("D9").Value = 45835
$ws.Range("D10").Value = 45835

# Move the active cell selection to F7 (matches the saved selection state)
$ws.Range("F7").Select()
